$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3) appended to the "Artfynd" sheet.
# Column letter -> column index reference (1-based):
#  A=1  B=2  C=3  D=4  E=5  F=6  G=7  H=8  I=9  J=10 K=11 L=12 M=13 N=14
#  O=15 P=16 Q=17 R=18 S=19 T=20 U=21 V=22 W=23 X=24 Y=25 Z=26 AA=27 AB=28
#  AC=29 AD=30 AE=31 AF=32 AG=33 AH=34 AI=35 AJ=36 AK=37 AL=38 AM=39 AN=40
#  AO=41 AP=42 AQ=43 AR=44 AS=45 AT=46 AU=47 AV=48 AW=49 AX=50 AY=51

$row = 3

$ws.Cells.Item($row, 1).Value  = 112215045
$ws.Cells.Item($row, 2).Value  = 56444
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"
$ws.Cells.Item($row, 4).Value  = "NT"
$ws.Cells.Item($row, 5).Value  = 102611
$ws.Cells.Item($row, 6).Value  = "Stenfalk"
$ws.Cells.Item($row, 7).Value  = "Falco columbarius"
$ws.Cells.Item($row, 8).Value  = "Linnaeus, 1758"

# Antal is stored as text "1" (matches source data, which keeps it as a string).
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value  = "1"

$ws.Cells.Item($row, 13).Value = "födosökande"

$ws.Cells.Item($row, 16).Value = "Norrtorp, Srm"
$ws.Cells.Item($row, 17).Value = 608521
$ws.Cells.Item($row, 18).Value = 6519076
$ws.Cells.Item($row, 19).Value = 3
$ws.Cells.Item($row, 20).Value = "Södermanland"
$ws.Cells.Item($row, 21).Value = "Nyköping"
$ws.Cells.Item($row, 22).Value = "Södermanland"
$ws.Cells.Item($row, 23).Value = "Stigtomta"

# Startdatum / Slutdatum are plain text strings ("2023-09-20"), not real
# Excel dates, so force text formatting before assigning the value.
$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2023-09-20"
$ws.Cells.Item($row, 26).Value = "13:01"
$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2023-09-20"
$ws.Cells.Item($row, 28).Value = "14:30"

$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false

$ws.Cells.Item($row, 49).Value = "Lillian Silfverduk"
$ws.Cells.Item($row, 50).Value = "Lillian Silfverduk"
